$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# employee_id, employee_name, department, absence_reason, absence_duration, absence_date, salary
$rows = @{
    2  = @(94008, "Theo da Cunha", "Juridico", "Viagem de negocios", 6, 45105, 6307.69)
    3  = @(65136, "Diego da Mota", "Operacoes", "Consulta medica", 7, 45094, 9617.549999999999)
    4  = @(73241, "Vinicius da Cruz", "Operacoes", "Problemas pessoais", 7, 45098, 6871.61)
    5  = @(1021, "Yuri Fogaça", "TI", "Viagem de negocios", 7, 45084, 9345.35)
    6  = @(31552, "Danilo Porto", "Recursos Humanos", "Consulta medica", 5, 45085, 6745.9)
    7  = @(37982, "Dom Sousa", "TI", "Consulta medica", 8, 45105, 7294.32)
    8  = @(72509, "Thiago Azevedo", "Marketing", "Problemas pessoais", 5, 45079, 7795.24)
    9  = @(94039, "Ravi Lucca Dias", "Marketing", "Problemas pessoais", 6, 45089, 8919)
    10 = @(11728, "Dra. Sophia Costa", "Vendas", "Viagem de negocios", 4, 45104, 9250.48)
    11 = @(49978, "Srta. Maria Sophia Silveira", "P&D", "Doenca", 6, 45104, 8535.34)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
    $ws.Cells.Item($r, 5).Value = $vals[4]
    $ws.Cells.Item($r, 6).Value = $vals[5]
    $ws.Cells.Item($r, 7).Value = $vals[6]
}
